$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main refs")

$ws.Range("A1").Value = "sequenceID"
$ws.Range("B1").Value = "species"
$ws.Range("C1").Value = "group"
$ws.Range("D1").Value = "subtype"
$ws.Range("F1").Value = "year"
$ws.Range("G1").Value = "country"
$ws.Range("H1").Value = "host"
$ws.Range("I1").Value = "host_common_name"

$ws.Range("H6").Select()
